# Remove the hyperlink around "Le Monde mathematical puzzle", keeping the
# run's text and character formatting intact (just un-linking it).
$d = $word.ActiveDocument

$h = $d.Hyperlinks(1)
$h.Delete()
